$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E7").Value = 15.4052
$ws.Range("B8").Value = 6.176599999999999
$ws.Range("B10").Value = 5.108899999999999
$ws.Range("B12").Value = 4.984300000000002
$ws.Range("E15").Value = 16.18719999999999
$ws.Range("B18").Value = 6.604500000000001
$ws.Range("E18").Value = 17.70050000000002
$ws.Range("E20").Value = 15.89809999999999
$ws.Range("E29").Value = 17.02830000000001
$ws.Range("E30").Value = 15.6101
$ws.Range("E31").Value = 16.2065
$ws.Range("B37").Value = 8.830400000000001
$ws.Range("E40").Value = 17.1711
$ws.Range("E50").Value = 16.2985
$ws.Range("B55").Value = 6.029199999999999
$ws.Range("B68").Value = 4.637199999999996
$ws.Range("E68").Value = 17.58150000000002
$ws.Range("E76").Value = 16.16319999999999
$ws.Range("B77").Value = 9.004500000000009
$ws.Range("B78").Value = 9.390300000000002
$ws.Range("B81").Value = 5.259100000000002
$ws.Range("B82").Value = 5.572099999999999
$ws.Range("E87").Value = 16.2977
$ws.Range("E88").Value = 16.3628
$ws.Range("E96").Value = 16.44499999999998
$ws.Range("E98").Value = 15.52
$ws.Range("E101").Value = 16.76170000000002
$ws.Range("E102").Value = 16.70189999999999
